# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (s_vals) values recomputed for rows 2-12 (column G)
$kValues = @{
    2  = 2
    3  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
